# Roll the retirement-planning templates forward by one year:
# drop the now-past "2024" row on both sheets so the tables start at 2025,
# and restore the originally-saved sheet activation / selection state.

$wb = $excel.ActiveWorkbook

$wsKim = $wb.Worksheets.Item("Kim")
$wsSam = $wb.Worksheets.Item("Sam")

# Remove the 2024 row from each sheet; Excel shifts everything below up by
# one row, which also updates the dimension and re-numbers all the years.
$wsKim.Rows.Item(2).Delete()
$wsSam.Rows.Item(2).Delete()

# Restore the saved view state: "Sam" had its selection set while it was the
# active sheet, then "Kim" was made active with its own selection - doing it
# in this order preserves both sheets' stored selections.
[void]$wsSam.Activate()
$wsSam.Range("A2:XFD2").Select() | Out-Null

[void]$wsKim.Activate()
$wsKim.Range("E21").Select() | Out-Null
